$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 112 ---
$ws.Range("E112").ClearFormats()
$ws.Range("F112").ClearFormats()
$ws.Cells.Item(112, 12).Value = 45910.65059689815
$ws.Cells.Item(112, 13).Value = 45910.65059133102

# --- Append new rows 113-128 ---
$newRows = @(
  @("Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED", "Letter-high", "(5)-NN_Classifier_GED", 0.2, 0.00004726797453703704, 0.00001169179398148148, 0.03111111111111111, 0.009831482319044507, 0.008878965089491406, 0.03111111111111111, 0, 45912.60330827547, 45912.60326099537, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED", "Letter-high", "(5)-NN_Classifier_GED", 0.2, 0.00003693153935185185, 0.000009869791666666667, 0.4711111111111111, 0.4643883570196614, 0.4972378117722602, 0.4711111111111111, 0, 45912.6047197338, 45912.60468278935, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED", "Letter-high", "(5)-NN_Classifier_GED", 0.2, 0.00003414002314814815, 0.000008256493055555557, 0.4911111111111111, 0.4829481352785313, 0.505079568830741, 0.4911111111111111, 0, 45912.60630778935, 45912.60627363426, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_Letter-high_with_(5)-NN_Classifier_GED", "Letter-high", "(5)-NN_Classifier_GED", 0.2, 0.00003405554398148148, 0.000008418587962962962, 0.5155555555555555, 0.507561490583848, 0.5272391631728707, 0.5155555555555555, 0, 45912.60910956019, 45912.60907547454, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000002130439814814815, 0.00000005905092592592592, 0.8421052631578947, 0.8362753036437246, 0.8520290625553782, 0.8421052631578947, 0, 45912.61058741898, 45912.6105871875, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000002159837962962963, 0.00000007855324074074074, 0.7368421052631579, 0.7319838056680164, 0.7314439946018894, 0.7368421052631579, 0, 45912.61173458333, 45912.61173435185, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.000000229537037037037, 0.00000007258101851851852, 0.9473684210526315, 0.9463967611336033, 0.951417004048583, 0.9473684210526315, 0.9285714285714286, 45912.61266902777, 45912.61266877314, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.000000183113425925926, 0.00000009755787037037038, 0.868421052631579, 0.8697626418988649, 0.8723886639676115, 0.868421052631579, 0.8589743589743589, 45912.612916875, 45912.61291667824, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.000000183113425925926, 0.00000009755787037037038, 0.9210526315789473, 0.9218575851393189, 0.9240485829959514, 0.9210526315789473, 0.9198717948717948, 45912.612916875, 45912.61291667824, "GEDLIB_Calculator", "Hyperparameter Tuning (grid)"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000002246180555555556, 0.00000005369212962962963, 0.7894736842105263, 0.7912679425837321, 0.7974347633790358, 0.7894736842105263, 0.7913043478260869, 45912.61332822917, 45912.61332799769, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000002086458333333334, 0.00000009187500000000001, 0.8421052631578947, 0.8385627530364372, 0.8403331561226298, 0.8421052631578947, 0.8061538461538462, 45912.61519488426, 45912.61519466435, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000005028587962962963, 0.0000001613773148148148, 0.8421052631578947, 0.8527327935222672, 0.9052631578947368, 0.8421052631578947, 0.896551724137931, 45912.61885483797, 45912.61885425926, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000004137615740740741, 0.0000001035185185185185, 0.9210526315789473, 0.9245386192754613, 0.9407894736842105, 0.9210526315789473, 0.9482758620689655, 45912.65649662037, 45912.65649607639, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000004137615740740741, 0.0000001035185185185185, 0.9210526315789473, 0.9245386192754613, 0.9407894736842105, 0.9210526315789473, 0.9482758620689655, 45912.65649662037, 45912.65649607639, "GEDLIB_Calculator", "Hyperparameter Tuning (grid)"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000002073958333333333, 0.00000007221064814814814, 0.8157894736842105, 0.8204112431616677, 0.8530701754385964, 0.8157894736842105, 0.8415384615384617, 45912.65889356482, 45912.65889335648, "GEDLIB_Calculator", "Simple Train-Test Split"),
  @("Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED", "MUTAG", "(5)-NN_Classifier_GED", 0.2, 0.0000002073958333333333, 0.00000007221064814814814, 0.868421052631579, 0.8710275321465806, 0.8845693779904306, 0.868421052631579, 0.8815384615384615, 45912.65889356853, 45912.65889335229, "GEDLIB_Calculator", "Hyperparameter Tuning (grid)")
)

$startRow = 113
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowdata = $newRows[$i]
    for ($c = 0; $c -lt $rowdata.Length; $c++) {
        $ws.Cells.Item($r, $c+1).Value = $rowdata[$c]
    }
    $ws.Cells.Item($r, 12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 13).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# Row 128 needs integer style (numFmtId 1, "0") on E and F, matching old row 112 style
$ws.Range("E128").NumberFormat = "0"
$ws.Range("F128").NumberFormat = "0"

Write-Host "Edit complete"